$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 27.142857
$ws.Range("I2").Value = 27.142857
$ws.Range("K2").Value = 27.142857
$ws.Range("M2").Value = 85.85714300000001
$ws.Range("H28").Value = 1213.75
$ws.Range("I28").Value = 605.9091
$ws.Range("J28").Value = 2551
$ws.Range("K28").Value = 605.9091
$ws.Range("L28").Value = 2551
$ws.Range("M28").Value = -120.9091
$ws.Range("N28").Value = -3521
$ws.Range("H111").Value = 10418905
$ws.Range("I111").Value = 31250832
$ws.Range("K111").Value = 93752496
$ws.Range("M111").Value = -93749429
$ws.Range("H132").Value = 2181.8125
$ws.Range("I132").Value = 1494.2727
$ws.Range("K132").Value = 4482.8181
$ws.Range("M132").Value = -1952.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5007969.5
$ws.Range("I32").Value = 5411170
$ws.Range("J32").Value = 35166.332
$ws.Range("K32").Value = 5411170
$ws.Range("L32").Value = 35166.332
$ws.Range("M32").Value = -5410883
$ws.Range("N32").Value = -35740.332
$ws.Range("H61").Value = 7539.727
$ws.Range("I61").Value = 1903.3334
$ws.Range("K61").Value = 1903.3334
$ws.Range("M61").Value = -1691.3334
$ws.Range("H86").Value = 70000
$ws.Range("J86").Value = 70000
$ws.Range("L86").Value = 70000
$ws.Range("N86").Value = -72372
$ws.Range("H89").Value = 70000
$ws.Range("J89").Value = 70000
$ws.Range("L89").Value = 210000
$ws.Range("N89").Value = -221856
$ws.Range("H110").Value = 41667996
$ws.Range("I110").Value = 1484.8334
$ws.Range("K110").Value = 1484.8334
$ws.Range("M110").Value = 560.1666
$ws.Range("H132").Value = 7093.771
$ws.Range("I132").Value = 5994.107
$ws.Range("J132").Value = 8633.299999999999
$ws.Range("K132").Value = 17982.321
$ws.Range("L132").Value = 25899.9
$ws.Range("M132").Value = -15452.321
$ws.Range("N132").Value = -30959.9
$ws.Range("H133").Value = 85696.57000000001
$ws.Range("J133").Value = 85696.57000000001
$ws.Range("L133").Value = 85696.57000000001
$ws.Range("N133").Value = -90756.57000000001
$ws.Range("H136").Value = 7539.727
$ws.Range("I136").Value = 1903.3334
$ws.Range("K136").Value = 5710.0002
$ws.Range("M136").Value = -3160.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 50730972
$ws.Range("I20").Value = 61407284
$ws.Range("J20").Value = 18499.25
$ws.Range("K20").Value = 61407284
$ws.Range("L20").Value = 18499.25
$ws.Range("M20").Value = -61407037
$ws.Range("N20").Value = -18993.25
$ws.Range("H86").Value = 26962.896
$ws.Range("I86").Value = 36340.395
$ws.Range("K86").Value = 36340.395
$ws.Range("M86").Value = -35217.395
$ws.Range("H89").Value = 26962.896
$ws.Range("I89").Value = 36340.395
$ws.Range("K89").Value = 181701.975
$ws.Range("M89").Value = -176085.975
$ws.Range("H105").Value = 6033.2104
$ws.Range("I105").Value = 6552.125
$ws.Range("K105").Value = 6552.125
$ws.Range("M105").Value = -4805.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6009.933
$ws.Range("I16").Value = 798.44446
$ws.Range("J16").Value = 8243.429
$ws.Range("K16").Value = 798.44446
$ws.Range("L16").Value = 8243.429
$ws.Range("M16").Value = -511.44446
$ws.Range("N16").Value = -8817.429
$ws.Range("H74").Value = 37666.332
$ws.Range("J74").Value = 37666.332
$ws.Range("L74").Value = 37666.332
$ws.Range("N74").Value = -39414.332
$ws.Range("H77").Value = 37666.332
$ws.Range("J77").Value = 37666.332
$ws.Range("L77").Value = 112998.996
$ws.Range("N77").Value = -121734.996
$ws.Range("H94").Value = 1732.3334
$ws.Range("I94").Value = 1955
$ws.Range("K94").Value = 1955
$ws.Range("M94").Value = -1504
$ws.Range("H99").Value = 6002
$ws.Range("I99").Value = 3751
$ws.Range("K99").Value = 3751
$ws.Range("M99").Value = -2253
$ws.Range("H113").Value = 6009.933
$ws.Range("I113").Value = 798.44446
$ws.Range("J113").Value = 8243.429
$ws.Range("K113").Value = 798.44446
$ws.Range("L113").Value = 8243.429
$ws.Range("M113").Value = 1371.55554
$ws.Range("N113").Value = -12583.429
$ws.Range("H126").Value = 6002
$ws.Range("I126").Value = 3751
$ws.Range("K126").Value = 11253
$ws.Range("M126").Value = -8783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 33534080
$ws.Range("J129").Value = 83833336
$ws.Range("L129").Value = 251500008
$ws.Range("N129").Value = -251510008
$ws.Range("H132").Value = 4637.659
$ws.Range("I132").Value = 2168.875
$ws.Range("J132").Value = 6048.393
$ws.Range("K132").Value = 19519.875
$ws.Range("L132").Value = 54435.537
$ws.Range("M132").Value = -16989.875
$ws.Range("N132").Value = -59495.537

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 4379.8076
$ws.Range("I102").Value = 4166.5454
$ws.Range("K102").Value = 4166.5454
$ws.Range("M102").Value = -2544.5454
$ws.Range("H107").Value = 657.5454999999999
$ws.Range("I107").Value = 437.9091
$ws.Range("J107").Value = 877.1818
$ws.Range("K107").Value = 437.9091
$ws.Range("L107").Value = 877.1818
$ws.Range("M107").Value = 1482.0909
$ws.Range("N107").Value = -4717.1818
$ws.Range("H122").Value = 65650.734
$ws.Range("I122").Value = 111934.27
$ws.Range("K122").Value = 335802.81
$ws.Range("M122").Value = -333352.81
$ws.Range("H126").Value = 4529.95
$ws.Range("I126").Value = 2527.0908
$ws.Range("J126").Value = 6977.8887
$ws.Range("K126").Value = 7581.2724
$ws.Range("L126").Value = 20933.6661
$ws.Range("M126").Value = -5111.2724
$ws.Range("N126").Value = -25873.6661
$ws.Range("H132").Value = 4823.7085
$ws.Range("I132").Value = 2134.7058
$ws.Range("J132").Value = 11354.143
$ws.Range("K132").Value = 6404.117400000001
$ws.Range("L132").Value = 34062.429
$ws.Range("M132").Value = -3874.117400000001
$ws.Range("N132").Value = -39122.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6117
$ws.Range("I7").Value = 3589.75
$ws.Range("K7").Value = 3589.75
$ws.Range("M7").Value = -3477.75
$ws.Range("H40").Value = 5880.5
$ws.Range("I40").Value = 3881.3333
$ws.Range("K40").Value = 3881.3333
$ws.Range("M40").Value = -3745.3333
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H68").Value = 5659.625
$ws.Range("I68").Value = 3359.1667
$ws.Range("K68").Value = 3359.1667
$ws.Range("M68").Value = -2610.1667
$ws.Range("H71").Value = 5659.625
$ws.Range("I71").Value = 3359.1667
$ws.Range("K71").Value = 16795.8335
$ws.Range("M71").Value = -13051.8335
$ws.Range("H122").Value = 4300.5454
$ws.Range("I122").Value = 2683.2727
$ws.Range("K122").Value = 8049.8181
$ws.Range("M122").Value = -5599.8181
$ws.Range("H126").Value = 6117
$ws.Range("I126").Value = 3589.75
$ws.Range("K126").Value = 10769.25
$ws.Range("M126").Value = -8299.25
$ws.Range("H132").Value = 7310.811
$ws.Range("I132").Value = 3692.3076
$ws.Range("K132").Value = 11076.9228
$ws.Range("M132").Value = -8546.9228
$ws.Range("H136").Value = 13231.585
$ws.Range("I136").Value = 3166.6667
$ws.Range("J136").Value = 21108.479
$ws.Range("K136").Value = 9500.000100000001
$ws.Range("L136").Value = 63325.437
$ws.Range("M136").Value = -6950.000100000001
$ws.Range("N136").Value = -68425.43700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21324
$ws.Range("I54").Value = 20000
$ws.Range("J54").Value = 21702.285
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 21702.285
$ws.Range("M54").Value = -19480
$ws.Range("N54").Value = -22742.285
$ws.Range("H122").Value = 3254.4688
$ws.Range("I122").Value = 2477.4167
$ws.Range("K122").Value = 7432.250100000001
$ws.Range("M122").Value = -4982.250100000001
$ws.Range("H126").Value = 1240
$ws.Range("I126").Value = 1400
$ws.Range("K126").Value = 4200
$ws.Range("M126").Value = -1730
$ws.Range("H132").Value = 5427.8486
$ws.Range("I132").Value = 7344.875
$ws.Range("J132").Value = 3623.5881
$ws.Range("K132").Value = 22034.625
$ws.Range("L132").Value = 10870.7643
$ws.Range("M132").Value = -19504.625
$ws.Range("N132").Value = -15930.7643
$ws.Range("H136").Value = 3952.186
$ws.Range("J136").Value = 6540.524
$ws.Range("L136").Value = 19621.572
$ws.Range("N136").Value = -24721.572
